# ---------------------------------------------------------------------------
# Operators test case update
#  - Inserts a new "Sheet1" worksheet (copy/rework of the operator sheet)
#    right before "usefullink".
#  - Updates wording + rich-text formatting for a couple of shared strings.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Fix the wording on the "operator" sheet (this also rewrites the shared
#    string table: removes the stale string and appends the corrected one).
# ---------------------------------------------------------------------------
$opSheet = $wb.Worksheets.Item("operator")
$opSheet.Range("C3").Value = "Technical Discussion (Errors, Software, Technical Materials)"

# ---------------------------------------------------------------------------
# 2. Insert the new "Sheet1" worksheet right before "usefullink".
# ---------------------------------------------------------------------------
$linkSheet = $wb.Worksheets.Item("usefullink")
$newSheet = $wb.Worksheets.Add($linkSheet)
$newSheet.Name = "Sheet1"

# Column widths (C, D, F)
$newSheet.Columns.Item(3).ColumnWidth = 52.0
$newSheet.Columns.Item(4).ColumnWidth = 25.333333333333332
$newSheet.Columns.Item(6).ColumnWidth = 19.833333333333332

# ---- Row 1 (header) --------------------------------------------------------
$newSheet.Range("A1").Value = "ID"
$newSheet.Range("B1").Value = "Person"
$newSheet.Range("C1").Value = "For"
$newSheet.Range("D1").Value = "Prefered Way to Connect"
$newSheet.Range("E1").Value = "Contact"
$newSheet.Range("F1").Value = "Timings"

# ---- Row 2 ------------------------------------------------------------------
$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = "Kiran"
$newSheet.Range("C2").Value = "Urgent Technical Help"
$newSheet.Range("D2").Value = "Whats App Only"
$newSheet.Range("E2").Value = 9552343698
$newSheet.Range("F2").Value = "07:00 AM to 10:00 PM Monday-Sunday"

# ---- Row 3 ------------------------------------------------------------------
$newSheet.Range("A3").Value = 2
$newSheet.Range("B3").Value = "Neelam"
$newSheet.Range("C3").Value = "Technical Discussion (Errors, Software, Technical Materials)"
$newSheet.Range("D3").Value = "Whats App Phone Call SMS eMail"
$newSheet.Range("E3").Value = 7066885937
$newSheet.Range("F3").Value = "09:00 AM to 06:00 PM Monday-Saturday"

# ---- Row 4 ------------------------------------------------------------------
$newSheet.Range("A4").Value = 3
$newSheet.Range("B4").Value = "Seema"
$newSheet.Range("C4").Value = "Administration (Fees, ID Card, Certificates, WhatsApp Group, Enquiry)"
$newSheet.Range("D4").Value = "Whats App Phone Call SMS eMail"
$newSheet.Range("E4").Value = 8888558802
$newSheet.Range("F4").Value = "09:00 AM to 06:00 PM Monday-Saturday"

# ---- Row 5 ------------------------------------------------------------------
$newSheet.Range("A5").Value = 4
$newSheet.Range("B5").Value = "Varsha"
$newSheet.Range("C5").Value = "Enquiry(Course Details, Fees, Enquiry)"
$newSheet.Range("D5").Value = "Whats App Phone Call SMS eMail"
$newSheet.Range("E5").Value = 8888809416
$newSheet.Range("F5").Value = "09:00 AM to 06:00 PM Monday to Friday and Sunday"

# ---- Row 6 ------------------------------------------------------------------
$newSheet.Range("A6").Value = 5
$newSheet.Range("B6").Value = "Darshit"
$newSheet.Range("C6").Value = "Technical Help"
$newSheet.Range("D6").Value = "Whats App Only"
$newSheet.Range("E6").Value = 8866888662
$newSheet.Range("F6").Value = "08:30 AM to 02:00 PM Saturday-Sunday"

$null = $newSheet.Range("H8").Select()
